# Applies the recorded edits to the "Artfynd" sheet.
# Net effect: row 4 and row 6 swap their species-observation data (with the
# coordinate values Q/R rounded to whole meters), rows 3 and 5 get their
# Q/R coordinates rounded to whole meters, and the now-unused Starttid/
# Sluttid ("00:00") cells on row 3 plus the comment cell that moved from
# row 4 to row 6 are cleared.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 3: round coordinates, drop the 00:00 start/end time cells -------
$ws.Range("Q3").Value = 511614
$ws.Range("R3").Value = 6733640
$ws.Range("Z3").ClearContents()
$ws.Range("AB3").ClearContents()

# --- Row 4: now holds the data previously shown on row 6 ------------------
$ws.Range("A4").Value = 112043031
$ws.Range("B4").Value = 90332
$ws.Range("E4").Value = 4769
$ws.Range("F4").Value = "Svavelriska"
$ws.Range("G4").Value = "Lactarius scrobiculatus"
$ws.Range("H4").Value = "(Scop.:Fr.) Fr."
$ws.Range("Q4").Value = 511625
$ws.Range("R4").Value = 6733616
$ws.Range("Z4").Value = "10:42"
$ws.Range("AB4").Value = "10:42"
$ws.Range("AC4").ClearContents()
$ws.Range("AX4").Value = "Evalena Sköld"

# --- Row 5: round coordinates ---------------------------------------------
$ws.Range("Q5").Value = 511628
$ws.Range("R5").Value = 6733623

# --- Row 6: now holds the data previously shown on row 4 ------------------
$ws.Range("A6").Value = 112042940
$ws.Range("B6").Value = 98535
$ws.Range("E6").Value = 222498
$ws.Range("F6").Value = "Blåsippa"
$ws.Range("G6").Value = "Hepatica nobilis"
$ws.Range("H6").Value = "Schreb."
$ws.Range("Q6").Value = 511611
$ws.Range("R6").Value = 6733626
$ws.Range("Z6").Value = "10:33"
$ws.Range("AB6").Value = "10:33"
$ws.Range("AC6").Value = "Fullt med blåsippsblad på denna sidan bäcken"
$ws.Range("AX6").Value = "Evalena Sköld, Åke Sköld"
